# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "34.960.26"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.846.52"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'232.94"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'41.60"
$ws.Range("E8").Value = "  +5.63%  "
$ws.Range("D9").Value = "'0.329"
$ws.Range("E9").Value = "  +3.08%  "
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("D11").Value = "'0.0982"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "2.114.37"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("D13").Value = "'11.43"
$ws.Range("E13").Value = "  +4.97%  "
$ws.Range("D14").Value = "1.854.39"
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("D17").Value = "34.957.84"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "'70.10"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").Value = "'240.88"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("E21").Value = "  +3.60%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").Value = "'172.74"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").Value = "'7.82"
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("D28").Value = "'0.123"
$ws.Range("D29").Value = "'1.68"
$ws.Range("E29").Value = "  +9.95%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("D32").Value = "'3.96"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").Value = "'3.93"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").Value = "'1.63"
$ws.Range("E34").Value = "  +23.29%  "
$ws.Range("E35").Value = "  +10.87%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.748"
$ws.Range("E36").Value = "  +9.44%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.24"
$ws.Range("E37").Value = "  +6.87%  "
$ws.Range("E38").Value = "  +11.69%  "
$ws.Range("D39").Value = "'89.86"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").Value = "1.349.51"
$ws.Range("E40").Value = "  +3.09%  "
$ws.Range("E41").Value = "  +2.91%  "
$ws.Range("D42").Value = "'14.64"
$ws.Range("E42").Value = "  +3.35%  "
$ws.Range("D43").Value = "'2.30"
$ws.Range("E43").Value = "  +4.71%  "
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("D46").Value = "'0.0532"
$ws.Range("D47").Value = "'6.35"
$ws.Range("E47").Value = "  +3.54%  "
$ws.Range("D48").Value = "2.031.98"
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("D49").Value = "'3.41"
$ws.Range("E49").Value = "  +15.93%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "'0.0670"
$ws.Range("E51").Value = "  -0.11%  "
